$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 269
$ws1.Range("F3").Value = 257
$ws1.Range("F4").Value = 277
$ws1.Range("F5").Value = 2887
$ws1.Range("F8").Value = 2243
$ws1.Range("F9").Value = 1430
$ws1.Range("F10").Value = 1430
$ws1.Range("F11").Value = 38
$ws1.Range("F12").Value = 449
$ws1.Range("F13").Value = 89
$ws1.Range("F14").Value = 2579
$ws1.Range("F16").Value = 1394
$ws1.Range("F17").Value = 4849
$ws1.Range("F19").Value = 5344
$ws1.Range("F20").Value = 5344
$ws1.Range("F21").Value = 1893
$ws1.Range("F22").Value = 2933
$ws1.Range("F23").Value = 3345
$ws1.Range("F24").Value = 190
$ws1.Range("F25").Value = 1611
$ws1.Range("F27").Value = 848
$ws1.Range("F28").Value = 129
$ws1.Range("F29").Value = 2
$ws1.Range("F30").Value = 322
$ws1.Range("F31").Value = 1043
$ws1.Range("F32").Value = 2088
$ws1.Range("F33").Value = 2
$ws1.Range("F34").Value = 128
$ws1.Range("F36").Value = 773
$ws1.Range("F38").Value = 369
$ws1.Range("F39").Value = 441

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 58
$ws2.Range("F15").Value = 13

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 58
$ws4.Range("F7").Value = 269
$ws4.Range("F9").Value = 277
$ws4.Range("F10").Value = 2887
$ws4.Range("F12").Value = 2243
$ws4.Range("F13").Value = 1430
$ws4.Range("F14").Value = 1430
$ws4.Range("F16").Value = 38
$ws4.Range("F17").Value = 449
$ws4.Range("F18").Value = 89
$ws4.Range("F20").Value = 2579
$ws4.Range("F21").Value = 1394
$ws4.Range("F26").Value = 4849
$ws4.Range("F28").Value = 5344
$ws4.Range("F29").Value = 5344
$ws4.Range("F30").Value = 1893
$ws4.Range("F31").Value = 2933
$ws4.Range("F32").Value = 3345
$ws4.Range("F33").Value = 190
$ws4.Range("F36").Value = 1611
$ws4.Range("F39").Value = 848
$ws4.Range("F40").Value = 129
$ws4.Range("F41").Value = 2
$ws4.Range("F42").Value = 322
$ws4.Range("F44").Value = 2088
$ws4.Range("F45").Value = 2
$ws4.Range("F46").Value = 128
$ws4.Range("F48").Value = 773
$ws4.Range("F50").Value = 369
$ws4.Range("F51").Value = 441
